$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated Price values are plain numerals (e.g. "312.49"); format those cells as
# Text first so Excel stores/keeps them as text like the source data, rather than
# coercing to a number (which would also silently drop significant trailing zeros).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.151.66"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "2.312.12"
$ws.Range("E3").Value = "  -2.41%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "312.49"
$ws.Range("E5").Value = "  -5.99%  "
$ws.Range("D6").Value = "105.70"
$ws.Range("E6").Value = "  +5.22%  "
$ws.Range("E7").Value = "  -1.94%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "0.609"
$ws.Range("E9").Value = "  -3.41%  "
$ws.Range("D10").Value = "40.15"
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("D11").Value = "0.0914"
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").Value = "8.28"
$ws.Range("E12").Value = "  -2.45%  "
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").Value = "0.978"
$ws.Range("E14").Value = "  -3.07%  "
$ws.Range("D15").Value = "15.56"
$ws.Range("E15").Value = "  -5.06%  "
$ws.Range("D16").Value = "2.657.34"
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("D17").Value = "2.298.04"
$ws.Range("E17").Value = "  -2.85%  "
$ws.Range("D18").Value = "42.148.44"
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("D19").Value = "7.70"
$ws.Range("E19").Value = "  -3.22%  "
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("D21").Value = "74.55"
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("D22").Value = "3.48"
$ws.Range("E22").Value = "  -8.09%  "
$ws.Range("D23").Value = "259.34"
$ws.Range("E23").Value = "  -3.58%  "
$ws.Range("D24").Value = "2.29"
$ws.Range("E24").Value = "  -1.38%  "
$ws.Range("D25").Value = "9.25"
$ws.Range("E25").Value = "  -7.56%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "11.00"
$ws.Range("E27").Value = "  -4.31%  "
$ws.Range("E28").Value = "  +3.34%  "
$ws.Range("D29").Value = "22.77"
$ws.Range("E29").Value = "  -1.70%  "
$ws.Range("D30").Value = "35.43"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").Value = "163.13"
$ws.Range("E32").Value = "  -7.57%  "
$ws.Range("E33").Value = "  -5.50%  "
$ws.Range("D34").Value = "5.85"
$ws.Range("E34").Value = "  -4.93%  "
$ws.Range("E35").Value = "  -2.72%  "
$ws.Range("D36").Value = "0.118"
$ws.Range("E36").Value = "  +11.88%  "
$ws.Range("D37").Value = "4.53"
$ws.Range("E37").Value = "  -1.66%  "
$ws.Range("D38").Value = "0.0352"
$ws.Range("E38").Value = "  -1.96%  "
$ws.Range("D39").Value = "2.76"
$ws.Range("E39").Value = "  -7.51%  "
$ws.Range("D40").Value = "3.63"
$ws.Range("E40").Value = "  -5.23%  "
$ws.Range("B41").Value = "BitcoinSV"
$ws.Range("C41").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D41").Value = "98.56"
$ws.Range("E41").Value = "  +8.52%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "1.46"
$ws.Range("E42").Value = "  -4.54%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").Value = "70.58"
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").Value = "0.230"
$ws.Range("E44").Value = "  -2.22%  "
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("D46").Value = "12.10"
$ws.Range("E46").Value = "  +1.47%  "
$ws.Range("D47").Value = "111.41"
$ws.Range("E47").Value = "  -6.02%  "
$ws.Range("D48").Value = "5.37"
$ws.Range("E48").Value = "  -2.17%  "
$ws.Range("D49").Value = "8.97"
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("D50").Value = "74.85"
$ws.Range("E50").Value = "  +6.03%  "
$ws.Range("D51").Value = "1.26"
$ws.Range("E51").Value = "  -0.97%  "
